$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting of the previous row down, then fill in the new values
$ws.Range("A3:C3").Copy()
$ws.Range("A4:C4").PasteSpecial(-4122)

$ws.Range("A4").Value = 44289
$ws.Range("B4").Value = "0.2.0"
$ws.Range("C4").Value = "All graphics from Cakedefi-Review.com integrated`nScripts are running on server and capture needed data`nFirst bugs are found and fixed"

$ws.Rows.Item(4).RowHeight = 45

$ws.Range("C5").Select()
